$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of cell B3 (previously contained "dd"), turning it into a blank cell.
$ws.Range("B3").ClearContents()

# Match the selection state recorded in the saved workbook.
$ws.Range("B3").Select()
